# This workbook holds a weekly price-reporting table (Acelga / Vega Monumental
# Concepción) where each week contributes a pair of rows ("Primera" / "Segunda").
# The update adds a new week's pair of rows at the top of the price-history
# block (rows 128-129, with a brand-new date) and pushes every following
# week's pair down by one position (2 rows), so the oldest week that used to
# sit at rows 179-180 now also reappears duplicated at the newly created
# rows 181-182.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: capture today's last two rows (179, 180) into the two brand-new
# rows (181, 182) before anything else is overwritten.
$ws.Range("A179:R179").Copy($ws.Range("A181:R181"))
$ws.Range("A180:R180").Copy($ws.Range("A182:R182"))

# Step 2: shift rows 130-180 down by 2 rows (row r takes the content that
# used to live at row r-2). Walk from the bottom up so the source row for
# each copy still holds its original, untouched content when it is read.
for ($r = 180; $r -ge 130; $r--) {
    $src = $r - 2
    $ws.Range("A" + $src + ":R" + $src).Copy($ws.Range("A" + $r + ":R" + $r))
}

# Step 3: rows 128-129 become the new week; only the date (column D) changes,
# every other value in those two rows stays as it already was.
$ws.Range("D128").Value2 = 44523
$ws.Range("D129").Value2 = 44523
